$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
